$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for changed rows ---
# D-column values are numeric-looking text; force text via leading apostrophe
# then reset the cell style back to Normal so no stray NumberFormat is applied.
$ws.Range("D2").Value = "'28.388.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.02%  "
$ws.Range("D3").Value = "'1.795.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'314.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.5445"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.69%  "
$ws.Range("D8").Value = "'0.3819"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.60%  "
$ws.Range("D9").Value = "'0.07567"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("D10").Value = "'42.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("E11").Value = "  +3.34%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("D13").Value = "'21.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.13%  "
$ws.Range("D14").Value = "'6.188"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").Value = "'7.390"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.13%  "
$ws.Range("D16").Value = "'1.794.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.52%  "
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("D18").Value = "'0.00001071"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("D19").Value = "'0.06445"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'17.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.33%  "
$ws.Range("D22").Value = "'5.957"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.25%  "
$ws.Range("D23").Value = "'28.431.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.02%  "
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").Value = "'2.127"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").Value = "'159.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.30%  "
$ws.Range("E27").Value = "  +2.36%  "
$ws.Range("D28").Value = "'2.387"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.23%  "
$ws.Range("D29").Value = "'2.006.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.65%  "
$ws.Range("D30").Value = "'123.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.90%  "
$ws.Range("D31").Value = "'1.127"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.75%  "
$ws.Range("D32").Value = "'0.1020"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.35%  "
$ws.Range("D33").Value = "'5.759"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.25%  "
$ws.Range("D37").Value = "'0.02324"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.37%  "
$ws.Range("D38").Value = "'5.152"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.42%  "
$ws.Range("D39").Value = "'8.746"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.96%  "
$ws.Range("D40").Value = "'11.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.31%  "
$ws.Range("D41").Value = "'0.6398"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.96%  "
$ws.Range("D42").Value = "'1.421"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "'0.9996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").Value = "'1.157"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.01%  "
$ws.Range("D45").Value = "'13.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.5987"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.68%  "
$ws.Range("D47").Value = "'3.667"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("D48").Value = "'126.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.91%  "
$ws.Range("D49").Value = "'2.006"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.15%  "
$ws.Range("D50").Value = "'1.152"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.30%  "
$ws.Range("D51").Value = "'0.06959"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.62%  "

# --- Rows 35/36 swapped (Algorand <-> Hedera) with updated price/volume ---
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.06805"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +13.96%  "

$ws.Range("B36").Value = "Algorand"
$ws.Range("C36").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D36").Value = "'0.2347"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +16.05%  "
